$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix tiny floating point rounding on G10 (recalculated diff = buy - sell) ---
$ws.Range("G10").Value = 0.07188224572437772

# --- Helper: write a literal text date (e.g. "2020-08-24") into a cell while
#     preserving the border/font/alignment formatting of a neighboring cell,
#     without letting Excel auto-convert the text into a date serial number. ---
function Set-DateText($cellRef, $sourceRef, $y, $m, $d) {
    $ws.Range($sourceRef).Copy($ws.Range($cellRef))
    $ws.Range($cellRef).Formula = "=TEXT(DATE($y,$m,$d),""yyyy-mm-dd"")"
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

# --- Row 11: 2020-08-24 ---
Set-DateText "A11" "A10" 2020 8 24
$ws.Range("B11").Value = 10974.5
$ws.Range("C11").Value = -0.2207990457459724
$ws.Range("D11").Value = 15507.3
$ws.Range("E11").Value = -0.2064589420680692
$ws.Range("F11").Value = 109541.8
$ws.Range("G11").Value = -0.3196833356208921
$ws.Range("H11").Value = 4700.299999999999
$ws.Range("I11").Value = -0.4572840532520466

# --- Row 12: 2020-08-25 ---
Set-DateText "A12" "A11" 2020 8 25
$ws.Range("B12").Value = 19049.6
$ws.Range("C12").Value = 0.7358057314684039
$ws.Range("D12").Value = 28799.3
$ws.Range("E12").Value = 0.8571446995930949
$ws.Range("F12").Value = 167733.5
$ws.Range("G12").Value = 0.5312282617229221
$ws.Range("H12").Value = 12142.2
$ws.Range("I12").Value = 1.583281918175436

# --- Row 13: 2020-08-26 ---
Set-DateText "A13" "A12" 2020 8 26
$ws.Range("B13").Value = 11468.3
$ws.Range("C13").Value = -0.3979768604065177
$ws.Range("D13").Value = 15837.7
$ws.Range("E13").Value = -0.4500664946717455
$ws.Range("F13").Value = 115022.1
$ws.Range("G13").Value = -0.3142568419546483
$ws.Range("H13").Value = 7992.799999999999
$ws.Range("I13").Value = -0.3417337879461713

# --- Row 14: 2020-08-27 ---
Set-DateText "A14" "A13" 2020 8 27
$ws.Range("B14").Value = 9505.5
$ws.Range("C14").Value = -0.1711500396745812
$ws.Range("D14").Value = 12331.3
$ws.Range("E14").Value = -0.2213957834786618
$ws.Range("F14").Value = 101931.2
$ws.Range("G14").Value = -0.1138120413381429
$ws.Range("H14").Value = 2607
$ws.Range("I14").Value = -0.6738314483034731

# --- Row 15: 2020-08-28 ---
Set-DateText "A15" "A14" 2020 8 28
$ws.Range("B15").Value = 32544.6
$ws.Range("C15").Value = 2.423765188575036
$ws.Range("D15").Value = 39411.8
$ws.Range("E15").Value = 2.196078272363823
$ws.Range("F15").Value = 321223.4
$ws.Range("G15").Value = 2.151374652706924
$ws.Range("H15").Value = 16661.2
$ws.Range("I15").Value = 5.390947449175298

Write-Host "Added rows 11-15 and corrected G10"
